$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "42.832.49"
$ws.Cells.Item(2, 5).Value = "  -0.08%  "
$ws.Cells.Item(3, 4).Value = "2.552.14"
$ws.Cells.Item(3, 5).Value = "  +0.69%  "
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "309.03"
$ws.Cells.Item(5, 5).Value = "  -3.02%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "101.93"
$ws.Cells.Item(6, 5).Value = "  +5.04%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.572"
$ws.Cells.Item(7, 5).Value = "  -0.34%  "
$ws.Cells.Item(8, 5).Value = "  +0.02%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.535"
$ws.Cells.Item(9, 5).Value = "  -0.42%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "36.29"
$ws.Cells.Item(10, 5).Value = "  +1.17%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0811"
$ws.Cells.Item(11, 5).Value = "  -1.05%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "7.42"
$ws.Cells.Item(12, 5).Value = "  -1.79%  "
$ws.Cells.Item(13, 5).Value = "  -0.62%  "
$ws.Cells.Item(14, 4).Value = "2.942.42"
$ws.Cells.Item(14, 5).Value = "  +0.66%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "15.94"
$ws.Cells.Item(15, 5).Value = "  +5.42%  "
$ws.Cells.Item(16, 4).Value = "2.542.77"
$ws.Cells.Item(16, 5).Value = "  +0.39%  "
$ws.Cells.Item(17, 5).Value = "  -1.21%  "
$ws.Cells.Item(18, 4).Value = "42.853.45"
$ws.Cells.Item(18, 5).Value = "  -0.11%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "6.78"
$ws.Cells.Item(19, 5).Value = "  -1.52%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "12.42"
$ws.Cells.Item(20, 5).Value = "  -2.34%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0958"
$ws.Cells.Item(21, 5).Value = "  -1.13%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "69.39"
$ws.Cells.Item(22, 5).Value = "  -0.44%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "247.76"
$ws.Cells.Item(23, 5).Value = "  -2.26%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.92"
$ws.Cells.Item(24, 5).Value = "  -1.82%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.08"
$ws.Cells.Item(25, 5).Value = "  +0.78%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "26.66"
$ws.Cells.Item(26, 5).Value = "  +0.82%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "40.81"
$ws.Cells.Item(28, 5).Value = "  -0.69%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.34"
$ws.Cells.Item(29, 5).Value = "  -3.06%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "10.16"
$ws.Cells.Item(30, 5).Value = "  -3.57%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "156.13"
$ws.Cells.Item(31, 5).Value = "  -0.96%  "
$ws.Cells.Item(32, 5).Value = "  -3.07%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0806"
$ws.Cells.Item(33, 5).Value = "  +1.65%  "
$ws.Cells.Item(34, 2).Value = "ARBITRUM"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.09"
$ws.Cells.Item(34, 5).Value = "  -2.92%  "
$ws.Cells.Item(35, 2).Value = "LidoDAOToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.30"
$ws.Cells.Item(35, 5).Value = "  -1.78%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.64"
$ws.Cells.Item(36, 5).Value = "  -2.52%  "
$ws.Cells.Item(37, 2).Value = "Celestia"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "18.48"
$ws.Cells.Item(37, 5).Value = "  -4.79%  "
$ws.Cells.Item(38, 2).Value = "ApeXProtocol"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.60"
$ws.Cells.Item(38, 5).Value = "  +5.80%  "
$ws.Cells.Item(39, 5).Value = "  -0.73%  "
$ws.Cells.Item(40, 5).Value = "  -0.75%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "4.29"
$ws.Cells.Item(41, 5).Value = "  +12.37%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "22.62"
$ws.Cells.Item(42, 5).Value = "  +3.59%  "
$ws.Cells.Item(43, 5).Value = "  -0.24%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0302"
$ws.Cells.Item(44, 5).Value = "  -1.08%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "3.29"
$ws.Cells.Item(45, 5).Value = "  +0.09%  "
$ws.Cells.Item(46, 4).Value = "1.987.54"
$ws.Cells.Item(46, 5).Value = "  -1.07%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "9.03"
$ws.Cells.Item(48, 4).Value = "2.795.22"
$ws.Cells.Item(48, 5).Value = "  +0.62%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "81.57"
$ws.Cells.Item(49, 5).Value = "  -3.19%  "
$ws.Cells.Item(50, 5).Value = "  +0.85%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "73.90"
$ws.Cells.Item(51, 5).Value = "  -1.83%  "
